$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Rows 2-4: average_area 3300 -> 220, with updated classification labels
$ws.Range("B2").Value = "All other industry"
$ws.Range("C2").Value = 220

$ws.Range("B3").Value = "Warehouses and storage"
$ws.Range("C3").Value = 220

$ws.Range("B4").Value = "Manufacturing and light industry"
$ws.Range("C4").Value = 220

# Rows 5-7: average_area 3300 -> 480, with updated classification labels
$ws.Range("B5").Value = "All other industry"
$ws.Range("C5").Value = 480

$ws.Range("B6").Value = "Warehouses and storage"
$ws.Range("C6").Value = 480

$ws.Range("B7").Value = "Manufacturing and light industry"
$ws.Range("C7").Value = 480
